# "Weighted stuff still isn't working :(" --------------------------------
# 1. Match Data: drop the duplicated placeholder stat-row data in F32:I33
#    (rows 32/33 keep only their "?" marker in column A).
# 2. Drive Team Data: refresh the (still-broken) weighted-average numbers;
#    several of them now blow up into #DIV/0! / #NUM! errors.
# 3. Leave the workbook with "Per Member Data" as the active/selected tab,
#    and restore "Match Data"'s selection to I32.

$wb = $excel.ActiveWorkbook

# --- Match Data --------------------------------------------------------
$wsMatch = $wb.Worksheets.Item("Match Data")
$wsMatch.Activate()
$wsMatch.Range("F32:I33").ClearContents() | Out-Null
$wsMatch.Range("I32").Select() | Out-Null

# --- Drive Team Data -----------------------------------------------------
$wsDrive = $wb.Worksheets.Item("Drive Team Data")
$wsDrive.Activate()

$wsDrive.Range("E2").Value = 63.528703703703705
$wsDrive.Range("F2").Value = 33.41111111111111
$wsDrive.Range("G2").Value = 20.605555555555554
$wsDrive.Range("H2").Value = "#NUM!"

$wsDrive.Range("E3").Value = "#DIV/0!"
$wsDrive.Range("F3").Value = 26.0
$wsDrive.Range("G3").Value = "#DIV/0!"
$wsDrive.Range("H3").Value = "#NUM!"

$wsDrive.Range("E4").Value = 70.23412698412699
$wsDrive.Range("F4").Value = 18.35
$wsDrive.Range("G4").Value = "#DIV/0!"
$wsDrive.Range("H4").Value = "#NUM!"

# --- Final active sheet/selection --------------------------------------
$wsPerMember = $wb.Worksheets.Item("Per Member Data")
$wsPerMember.Activate()
